# Update benchmark: 2026-01-27 06:43:59 UTC
# Clears stale benchmark values from cells that no longer apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @(
    "E2",
    "C3", "F3",
    "C4", "F4",
    "C5", "F5",
    "C6",
    "C8", "F8",
    "C9", "F9",
    "C10", "F10",
    "C11",
    "C12",
    "C13", "F13",
    "C14", "F14"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
